$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of existing headers (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and J (IF), rows 2..45
$iValues = @(8,7,8,8,9,6,7,7,8,9,8,8,7,6,1,8,6,8,1,5,10,6,9,9,7,7,7,7,8,5,5,8,8,8,8,6,6,7,7,8,7,5,6,7)
$jValues = @(8,7,8,8,9,6,7,9,8,9,8,8,8,6,2,8,6,9,2,5,10,6,9,9,7,7,7,7,8,6,5,8,8,8,8,6,6,7,7,8,7,5,6,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
